$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing cell contents (keep formatting/styles) so the shared-string table rebuilds
# in the exact order the new data requires.
$ws.Cells.ClearContents()

# --- Row 1: headers (order unchanged) ---
$ws.Range("B1").Value = "state"
$ws.Range("C1").Value = "precision"
$ws.Range("D1").Value = "recall"
$ws.Range("E1").Value = "fmeasure"
$ws.Range("F1").Value = "final_precision"
$ws.Range("G1").Value = "final_recall"
$ws.Range("H1").Value = "final_fmeasure"
$ws.Range("I1").Value = "operation"
$ws.Range("J1").Value = "delta_related"
$ws.Range("K1").Value = "delta_visited"
$ws.Range("L1").Value = "delta_selected"
$ws.Range("M1").Value = "related"
$ws.Range("N1").Value = "visited"
$ws.Range("O1").Value = "final_selected"
$ws.Range("P1").Value = "accumulated_precision"
$ws.Range("Q1").Value = "accumulated_recall"
$ws.Range("R1").Value = "accumulated_fmeasure"
$ws.Range("S1").Value = "final_accumulated_precision"
$ws.Range("T1").Value = "final_accumulated_recall"
$ws.Range("U1").Value = "final_accumulated_fmeasure"

# --- Column A: row index numbers ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

# --- Column B: state names, written top-to-bottom in the NEW order ---
$ws.Range("B2").Value = "webofscience"
$ws.Range("B3").Value = "springer"
$ws.Range("B4").Value = "scopus"
$ws.Range("B5").Value = "sciencedirect"
$ws.Range("B6").Value = "ieee"
$ws.Range("B7").Value = "elcompendex"
$ws.Range("B8").Value = "acm"
$ws.Range("B9").Value = "s0"

# --- Column I: operation strings ("-" for rows 2-8, "union" for row 9) ---
$ws.Range("I2").Value = "-"
$ws.Range("I3").Value = "-"
$ws.Range("I4").Value = "-"
$ws.Range("I5").Value = "-"
$ws.Range("I6").Value = "-"
$ws.Range("I7").Value = "-"
$ws.Range("I8").Value = "-"
$ws.Range("I9").Value = "union"

# --- Remaining numeric columns ---
# Row 2
$ws.Range("C2").Value = 50
$ws.Range("D2").Value = 9.8
$ws.Range("E2").Value = 16.39
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 50
$ws.Range("Q2").Value = 9.8
$ws.Range("R2").Value = 16.39
$ws.Range("S2").Value = 0.4
$ws.Range("T2").Value = 0.1333333333333333
$ws.Range("U2").Value = 0.2

# Row 3
$ws.Range("C3").Value = 1.42
$ws.Range("D3").Value = 3.92
$ws.Range("E3").Value = 2.08
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 141
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 141
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1.42
$ws.Range("Q3").Value = 3.92
$ws.Range("R3").Value = 2.08
$ws.Range("S3").Value = 0.007092198581560284
$ws.Range("T3").Value = 0.03333333333333333
$ws.Range("U3").Value = 0.01169590643274854

# Row 4
$ws.Range("C4").Value = 46.67
$ws.Range("D4").Value = 13.73
$ws.Range("E4").Value = 21.21
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 15
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 7
$ws.Range("N4").Value = 15
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = 46.67
$ws.Range("Q4").Value = 13.73
$ws.Range("R4").Value = 21.21
$ws.Range("S4").Value = 0.4
$ws.Range("T4").Value = 0.2
$ws.Range("U4").Value = 0.2666666666666667

# Row 5
$ws.Range("C5").Value = 0.51
$ws.Range("D5").Value = 1.96
$ws.Range("E5").Value = 0.8099999999999999
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 195
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 195
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 0.51
$ws.Range("Q5").Value = 1.96
$ws.Range("R5").Value = 0.8099999999999999
$ws.Range("S5").Value = 0.005128205128205128
$ws.Range("T5").Value = 0.03333333333333333
$ws.Range("U5").Value = 0.008888888888888889

# Row 6
$ws.Range("C6").Value = 13.95
$ws.Range("D6").Value = 11.76
$ws.Range("E6").Value = 12.77
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 43
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 43
$ws.Range("O6").Value = 4
$ws.Range("P6").Value = 13.95
$ws.Range("Q6").Value = 11.76
$ws.Range("R6").Value = 12.77
$ws.Range("S6").Value = 0.09302325581395349
$ws.Range("T6").Value = 0.1333333333333333
$ws.Range("U6").Value = 0.1095890410958904

# Row 7
$ws.Range("C7").Value = 38.46
$ws.Range("D7").Value = 9.8
$ws.Range("E7").Value = 15.62
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 13
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 5
$ws.Range("P7").Value = 38.46
$ws.Range("Q7").Value = 9.8
$ws.Range("R7").Value = 15.62
$ws.Range("S7").Value = 0.3846153846153846
$ws.Range("T7").Value = 0.1666666666666667
$ws.Range("U7").Value = 0.2325581395348837

# Row 8
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 9.8
$ws.Range("E8").Value = 6.619999999999999
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 100
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = 5
$ws.Range("Q8").Value = 9.8
$ws.Range("R8").Value = 6.619999999999999
$ws.Range("S8").Value = 0.03
$ws.Range("T8").Value = 0.1
$ws.Range("U8").Value = 0.04615384615384615

# Row 9
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 22
$ws.Range("N9").Value = 497
$ws.Range("O9").Value = 15
$ws.Range("P9").Value = 4.43
$ws.Range("Q9").Value = 43.14
$ws.Range("R9").Value = 8.03
$ws.Range("S9").Value = 0.03018108651911469
$ws.Range("T9").Value = 0.5
$ws.Range("U9").Value = 0.05692599620493359
